$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The J2:K10 columns used to hold "=TRUE()" boolean formulas (Yes/No indicator).
# They are converted into plain "Yes" text values, formatted with a custom
# number format that still displays as TRUE/TRUE/FALSE.
$boolRange = $ws.Range("J2:K10")
$boolRange.NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$boolRange.Value = "Yes"

# L1 header cell picks up the same "wrap-capable" style used elsewhere in row 1.
$ws.Range("L1").Style = $ws.Range("A1").Style

# Selection moves from L2 to the newly updated K2:K10 range.
$ws.Range("K2:K10").Select()
